$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu results for the 380 kV case (rows 2-25 => res_bus index 0-23).
# Column layout: B..F = bus voltages for columns 0-4, G (bus 5, slack) stays 1,
# H is unused, I..M = bus voltages for columns 7-11.
$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13 }

$newValues = @{
    2 = @{ "B" = 1.02; "C" = 1.037579597792125; "D" = 1.04053330074585; "E" = 1.045203871433052; "F" = 1.053372110736915; "I" = 1.038082771872998; "J" = 1.042681699862228; "K" = 1.043315379284932; "L" = 1.047972769223769; "M" = 1.056118261612939 }
    3 = @{ "B" = 1.02; "C" = 1.038466054398504; "D" = 1.041199342010531; "E" = 1.046026717697069; "F" = 1.054356577663086; "I" = 1.038279962851485; "J" = 1.043212910844094; "K" = 1.043792243320464; "L" = 1.048606978714172; "M" = 1.056915320911727 }
    4 = @{ "B" = 1.02; "C" = 1.039040177727295; "D" = 1.041630689919628; "E" = 1.046560055103997; "F" = 1.054994846577108; "I" = 1.038406602487467; "J" = 1.043556519681174; "K" = 1.044100498083387; "L" = 1.049017604598287; "M" = 1.05743170564266 }
    5 = @{ "B" = 1.02; "C" = 1.0392816639372; "D" = 1.0418121167198; "E" = 1.046784484269654; "F" = 1.055263473244372; "I" = 1.038459612283473; "J" = 1.043700942959238; "K" = 1.044230013493026; "L" = 1.049190290323008; "M" = 1.057648944199868 }
    6 = @{ "B" = 1.02; "C" = 1.039322217812475; "D" = 1.041842584201065; "E" = 1.046822179425273; "F" = 1.055308594298051; "I" = 1.03846849938534; "J" = 1.043725190479064; "K" = 1.044251755275923; "L" = 1.04921928841677; "M" = 1.057685428289927 }
    7 = @{ "B" = 1.02; "C" = 1.039043403989255; "D" = 1.041633113809042; "E" = 1.046563053098084; "F" = 1.05499843480915; "I" = 1.03840731170952; "J" = 1.04355844959086; "K" = 1.044102228970307; "L" = 1.049019911804664; "M" = 1.057434607805585 }
    8 = @{ "B" = 1.02; "C" = 1.037879070119657; "D" = 1.040758313994126; "E" = 1.045481768723316; "F" = 1.05370455557189; "I" = 1.03814961108862; "J" = 1.042861249322774; "K" = 1.043476600850452; "L" = 1.048187050798431; "M" = 1.056387499416241 }
    9 = @{ "B" = 1.02; "C" = 1.035831456565464; "D" = 1.039219743367256; "E" = 1.043583375320005; "F" = 1.051434240041506; "I" = 1.037688214077042; "J" = 1.04163182166446; "K" = 1.042371853577857; "L" = 1.046721411990372; "M" = 1.054547277931169 }
    10 = @{ "B" = 1.02; "C" = 1.034469209463584; "D" = 1.038196097926452; "E" = 1.042322556366991; "F" = 1.049927286629879; "I" = 1.037375747318874; "J" = 1.040811682633901; "K" = 1.041633868408017; "L" = 1.045745716251056; "M" = 1.053323851674419 }
    11 = @{ "B" = 1.02; "C" = 1.033880028501522; "D" = 1.037753358224991; "E" = 1.04177776008618; "F" = 1.049276341558851; "I" = 1.03723929857591; "J" = 1.040456443611392; "K" = 1.041313973936863; "L" = 1.04532357585312; "M" = 1.052794915357921 }
    12 = @{ "B" = 1.02; "C" = 1.033661284047918; "D" = 1.037588982503987; "E" = 1.041575572309051; "F" = 1.049034789847472; "I" = 1.037188443469678; "J" = 1.040324476069615; "K" = 1.041195100609278; "L" = 1.045166826760283; "M" = 1.052598568719622 }
    13 = @{ "B" = 1.02; "C" = 1.033708200767651; "D" = 1.037624238121168; "E" = 1.041618934383062; "F" = 1.049086592702484; "I" = 1.037199359834959; "J" = 1.040352784259686; "K" = 1.041220601598951; "L" = 1.045200447570799; "M" = 1.052640680113122 }
    14 = @{ "B" = 1.02; "C" = 1.033861944896183; "D" = 1.037739769286231; "E" = 1.041761043626646; "F" = 1.049256369949453; "I" = 1.037235098381618; "J" = 1.040445535452341; "K" = 1.041304148849763; "L" = 1.045310617843163; "M" = 1.052778682748915 }
    15 = @{ "B" = 1.02; "C" = 1.033956685492684; "D" = 1.037810962168662; "E" = 1.041848624900794; "F" = 1.049361006934063; "I" = 1.037257095310113; "J" = 1.040502680442082; "K" = 1.041355618439256; "L" = 1.045378504396796; "M" = 1.052863727181154 }
    16 = @{ "B" = 1.02; "C" = 1.034508325920529; "D" = 1.038225491887015; "E" = 1.042358736992529; "F" = 1.049970521039989; "I" = 1.037384778823947; "J" = 1.04083525638118; "K" = 1.041655091638992; "L" = 1.045773739634916; "M" = 1.053358972730337 }
    17 = @{ "B" = 1.02; "C" = 1.034854538613015; "D" = 1.038485651878725; "E" = 1.042679024912835; "F" = 1.050353276081409; "I" = 1.037464564230409; "J" = 1.041043842884809; "K" = 1.041842852477415; "L" = 1.046021752761304; "M" = 1.053669846429195 }
    18 = @{ "B" = 1.02; "C" = 1.035056544242249; "D" = 1.038637447471557; "E" = 1.04286595396874; "F" = 1.050576682550145; "I" = 1.037510990847459; "J" = 1.041165496847957; "K" = 1.041952337123256; "L" = 1.046166447566498; "M" = 1.053851252373479 }
    19 = @{ "B" = 1.02; "C" = 1.035125434033692; "D" = 1.03868921405376; "E" = 1.042929710646993; "F" = 1.050652884105597; "I" = 1.037526802302213; "J" = 1.041206975804127; "K" = 1.041989662914463; "L" = 1.04621579030669; "M" = 1.05391312039872 }
    20 = @{ "B" = 1.02; "C" = 1.03481738647072; "D" = 1.038457734128135; "E" = 1.042644649606257; "F" = 1.050312194361989; "I" = 1.037456015473939; "J" = 1.041021464658536; "K" = 1.04182271092946; "L" = 1.04599513988358; "M" = 1.053636484486098 }
    21 = @{ "B" = 1.02; "C" = 1.033816668212491; "D" = 1.037705746085228; "E" = 1.041719191185708; "F" = 1.049206368191966; "I" = 1.037224579012767; "J" = 1.040418222969479; "K" = 1.041279547645772; "L" = 1.045278173980754; "M" = 1.052738040982149 }
    22 = @{ "B" = 1.02; "C" = 1.033188076892932; "D" = 1.037233390621379; "E" = 1.041138325173445; "F" = 1.048512470718222; "I" = 1.037078071066837; "J" = 1.040038849224013; "K" = 1.040937748988941; "L" = 1.044827694264821; "M" = 1.052173871189809 }
    23 = @{ "B" = 1.02; "C" = 1.033521247748358; "D" = 1.037483752081897; "E" = 1.041446157321926; "F" = 1.048880187638313; "I" = 1.037155831784493; "J" = 1.040239970758251; "K" = 1.041118970075021; "L" = 1.045066472747763; "M" = 1.052472879822405 }
    24 = @{ "B" = 1.02; "C" = 1.034834173721989; "D" = 1.038470348807724; "E" = 1.04266018198629; "F" = 1.050330756953386; "I" = 1.037459878631615; "J" = 1.041031576449884; "K" = 1.041831812130317; "L" = 1.046007164995828; "M" = 1.0536515590695 }
    25 = @{ "B" = 1.02; "C" = 1.036360320897136; "D" = 1.039617142982558; "E" = 1.044073320333561; "F" = 1.052020016198537; "I" = 1.037808357087609; "J" = 1.041949754742331; "K" = 1.042657724194693; "L" = 1.047100073920643; "M" = 1.055022428136482 }
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value = $rowValues[$col]
    }
}
